# Update "想去人数" (interest counts) in column F across all sheets
# as published in the regenerated gh-pages data at commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1740
$ws.Range("F4").Value = 144
$ws.Range("F8").Value = 1114
$ws.Range("F9").Value = 296
$ws.Range("F11").Value = 859
$ws.Range("F12").Value = 650
$ws.Range("F13").Value = 174
$ws.Range("F15").Value = 139
$ws.Range("F18").Value = 2854
$ws.Range("F19").Value = 2598
$ws.Range("F20").Value = 523
$ws.Range("F26").Value = 2436
$ws.Range("F29").Value = 15
$ws.Range("F31").Value = 270
$ws.Range("F32").Value = 1055

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1066
$ws.Range("F5").Value = 1066
$ws.Range("F16").Value = 98
$ws.Range("F25").Value = 303
$ws.Range("F27").Value = 3833
$ws.Range("F32").Value = 35
$ws.Range("F34").Value = 143
$ws.Range("F36").Value = 6

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1772
$ws.Range("F5").Value = 2423
$ws.Range("F6").Value = 1007
$ws.Range("F9").Value = 1268

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1772
$ws.Range("F3").Value = 2423
$ws.Range("F4").Value = 1740
$ws.Range("F6").Value = 1007
$ws.Range("F7").Value = 1268
$ws.Range("F10").Value = 144
$ws.Range("F15").Value = 1114
$ws.Range("F16").Value = 296
$ws.Range("F17").Value = 650
$ws.Range("F18").Value = 1066
$ws.Range("F19").Value = 174
$ws.Range("F23").Value = 2854
$ws.Range("F24").Value = 2598
$ws.Range("F31").Value = 2436
$ws.Range("F36").Value = 15
$ws.Range("F37").Value = 98
$ws.Range("F38").Value = 270
$ws.Range("F42").Value = 303
$ws.Range("F44").Value = 1055
$ws.Range("F46").Value = 35
$ws.Range("F48").Value = 143

$wb.Save()
